$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 17 (push old rows 17.. down), then copy the formatting
# of the last existing data row (16) onto the new A17:A19 cells so they match the
# existing "index" column style instead of Excel's auto-extended (borderless) style.
$ws.Rows.Item(17).Resize(3).Insert()
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# Rewrite rows 10-19 (A,B labels + C:M averaged-intensity values) to reflect the new
# method ordering (Gaussian-Quadrature + 3 new Spiral-* schemes) per the latest run.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.001031114751944
$ws.Range("D10").Value = 0.9872471667522974
$ws.Range("E10").Value = 1.000508245343851
$ws.Range("F10").Value = 1.001031114751944
$ws.Range("G10").Value = 0.9902760876083408
$ws.Range("H10").Value = 1.002854316891267
$ws.Range("I10").Value = 0.9994117647058823
$ws.Range("J10").Value = 0.9872471667522974
$ws.Range("K10").Value = 0.9938777060480739
$ws.Range("L10").Value = 0.9974544104000088
$ws.Range("M10").Value = 0.9968881160089302

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9985927653537351
$ws.Range("D11").Value = 0.9857461190329304
$ws.Range("E11").Value = 1.002659292154276
$ws.Range("F11").Value = 0.9985927653537351
$ws.Range("G11").Value = 0.9900964061099661
$ws.Range("H11").Value = 1.009765281260677
$ws.Range("I11").Value = 1.001505840299995
$ws.Range("J11").Value = 0.9857461190329304
$ws.Range("K11").Value = 0.994202705593603
$ws.Range("L11").Value = 0.9963977354736689
$ws.Range("M11").Value = 0.9980609507019298

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9987279126857799
$ws.Range("D12").Value = 0.9851855712761802
$ws.Range("E12").Value = 1.002803202256619
$ws.Range("F12").Value = 0.9987279126857799
$ws.Range("G12").Value = 0.9897465978872384
$ws.Range("H12").Value = 1.010162071374739
$ws.Range("I12").Value = 1.001643922344458
$ws.Range("J12").Value = 0.9851855712761802
$ws.Range("K12").Value = 0.9939943867663994
$ws.Range("L12").Value = 0.9963611497260897
$ws.Range("M12").Value = 0.9980448796375022

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9986261406762679
$ws.Range("D13").Value = 0.9856156108243287
$ws.Range("E13").Value = 1.002692941698999
$ws.Range("F13").Value = 0.9986261406762679
$ws.Range("G13").Value = 0.9899938880961293
$ws.Range("H13").Value = 1.009924647781642
$ws.Range("I13").Value = 1.001535060893574
$ws.Range("J13").Value = 0.9856156108243287
$ws.Range("K13").Value = 0.994154276261664
$ws.Range("L13").Value = 0.9963902084689659
$ws.Range("M13").Value = 0.9980647149951568

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.011868000000001
$ws.Range("D14").Value = 0.9358720000000013
$ws.Range("E14").Value = 1.014780000000001
$ws.Range("F14").Value = 1.011868000000001
$ws.Range("G14").Value = 0.9666959999999998
$ws.Range("H14").Value = 1.025792000000001
$ws.Range("I14").Value = 1.013532000000002
$ws.Range("J14").Value = 0.9358720000000013
$ws.Range("K14").Value = 0.975326000000001
$ws.Range("L14").Value = 0.9935970000000007
$ws.Range("M14").Value = 0.9947566666666674

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.01
$ws.Range("D15").Value = 0.95
$ws.Range("E15").Value = 1.01
$ws.Range("F15").Value = 1.01
$ws.Range("G15").Value = 0.9793124999999996
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1.01
$ws.Range("J15").Value = 0.95
$ws.Range("K15").Value = 0.98
$ws.Range("L15").Value = 0.9949999999999999
$ws.Range("M15").Value = 0.9932187499999999

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.005335203635197
$ws.Range("D16").Value = 0.9700070047744059
$ws.Range("E16").Value = 1.004678046105596
$ws.Range("F16").Value = 1.005335203635197
$ws.Range("G16").Value = 0.986820929433605
$ws.Range("H16").Value = 1.000626159411196
$ws.Range("I16").Value = 1.004776022220794
$ws.Range("J16").Value = 0.9700070047744059
$ws.Range("K16").Value = 0.9873425254400008
$ws.Range("L16").Value = 0.996338864537599
$ws.Range("M16").Value = 0.9953738942634658

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9978853896048567
$ws.Range("D17").Value = 0.9977727088113082
$ws.Range("E17").Value = 0.9979790335664273
$ws.Range("F17").Value = 0.9978853896048567
$ws.Range("G17").Value = 0.9979669216176827
$ws.Range("H17").Value = 0.9983359572995982
$ws.Range("I17").Value = 0.9979397601024624
$ws.Range("J17").Value = 0.9977727088113082
$ws.Range("K17").Value = 0.9978758711888678
$ws.Range("L17").Value = 0.9978806303968621
$ws.Range("M17").Value = 0.9979799618337227

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9983287630965081
$ws.Range("D18").Value = 0.9991030917230529
$ws.Range("E18").Value = 0.9976735979817106
$ws.Range("F18").Value = 0.9983287630965081
$ws.Range("G18").Value = 0.9999229021726924
$ws.Range("H18").Value = 0.9960141506691166
$ws.Range("I18").Value = 0.9974210105043452
$ws.Range("J18").Value = 0.9991030917230529
$ws.Range("K18").Value = 0.9983883448523818
$ws.Range("L18").Value = 0.998358553974445
$ws.Range("M18").Value = 0.9980772526912376

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9967851908990601
$ws.Range("D19").Value = 1.002657155842219
$ws.Range("E19").Value = 0.9967265091471618
$ws.Range("F19").Value = 0.9967851908990601
$ws.Range("G19").Value = 1.000829627201395
$ws.Range("H19").Value = 0.994906406510452
$ws.Range("I19").Value = 0.9965038206260239
$ws.Range("J19").Value = 1.002657155842219
$ws.Range("K19").Value = 0.9996918324946902
$ws.Range("L19").Value = 0.9982385116968752
$ws.Range("M19").Value = 0.998068118371052

Write-Output "applied averaged-intensities update for spiral schemes"
